$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the bold/bordered/centered style
# already used by the other header cells (e.g. G1 "sum").
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for the "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0

$excel.CutCopyMode = 0
